$d = $word.ActiveDocument

# The paragraph currently reads "Version 1." as:
#   [proofErr spellStart] "Version" [proofErr spellEnd] " 1." [bookmarkStart _GoBack][bookmarkEnd]
# and must become "Version 2." laid out as:
#   [proofErr spellStart] "Versi" "on" [proofErr spellEnd] " 2" [bookmarkStart _GoBack][bookmarkEnd] "."
#
# i.e. besides the digit change, the word "Version" gets split into two runs
# ("Versi"/"on") and the trailing "." moves to its own run positioned after the
# _GoBack bookmark instead of before it.

# 1) Force a run break inside "Version" at offset 5 (between "Versi" and "on")
#    without altering any visible formatting. Toggling a character property on
#    and back off on just the "Versi" sub-range is enough to make the engine
#    keep that run boundary instead of re-merging the word into a single run.
$splitPoint = $d.Range(0, 5)
$splitPoint.Font.Bold = $true
$splitPoint.Font.Bold = $false

# 2) The digit lives inside the " 1." run (character offsets 7-10, i.e. the
#    space, "1" and "."). Replace it with " 2", dropping the trailing period
#    here - it gets reinserted after the bookmark in the next step so that it
#    ends up on the correct side of _GoBack.
$numRange = $d.Range(7, 10)
$numRange.Text = " 2"

# 3) Re-append the "." so it lands AFTER the _GoBack bookmark. The bookmark
#    sits at the end of the paragraph's visible text, and inserting through
#    the whole-document range collapsed to its very end places new content
#    after the bookmark (matching the target markup), unlike inserting via a
#    range collapsed at the same character offset, which lands before it.
$tail = $d.Content
$tail.Collapse(0)
$tail.InsertAfter(".")
